# Commit: Wed, Jul 22, 2020  7:05:33 PM
#
# 1) Slide 6 ("SOURCES OF FINANCE") table switches from the custom
#    "Table_0" table style to the built-in table-style id
#    {83025033-CB8E-4329-BB2B-4A75441CDE80}.
# 2) The deck's active theme ("Integral", ppt/theme/theme2.xml) has its
#    12 theme colors swapped for the stock "Office" color scheme that
#    used to live in the (otherwise unused) ppt/theme/theme1.xml part.

$p = $ppt.ActivePresentation

# --- 1) Re-style the table on slide 6 -------------------------------
$s6 = $p.Slides.Item(6)
for ($i = 1; $i -le $s6.Shapes.Count; $i++) {
    $shp = $s6.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{83025033-CB8E-4329-BB2B-4A75441CDE80}", $false)
    }
}

# --- 2) Swap the active theme's color scheme to the "Office" palette -
$s1 = $p.Slides.Item(1)
$tcs = $s1.ThemeColorScheme

# order matches OOXML a:clrScheme children:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeColors = @(
    0,        # dk1      000000
    16777215, # lt1      FFFFFF
    6968388,  # dk2      44546A
    15132391, # lt2      E7E6E6
    13998939, # accent1  5B9BD5
    3243501,  # accent2  ED7D31
    10855845, # accent3  A5A5A5
    49407,    # accent4  FFC000
    12874308, # accent5  4472C4
    4697456,  # accent6  70AD47
    12673797, # hlink    0563C1
    7491477   # folHlink 954F72
)

for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Item($i).RGB = $officeColors[$i - 1]
}
